$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 4
    3  = -2
    4  = -2
    5  = 1
    6  = -3
    7  = -4
    8  = -4
    9  = 1
    10 = 0
    11 = -4
    13 = -4
    14 = -4
    15 = 5
    16 = -1
    17 = 2
    18 = -1
    19 = 3
    20 = -1
    21 = 2
    22 = 3
    23 = -4
    24 = 3
    26 = 4
    27 = -4
    28 = -1
    29 = -1
    30 = -1
    32 = 0
    33 = 2
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
